# cv121062a.xlsx - "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet had three label-only rows that carried no data (placeholder
# section headers / footer note):
#   row 5  -> "situação do domicílio"                (section header, blank)
#   row 8  -> "grandes regiões e unidades da federação" (section header, blank)
#   row 41 -> "fonte: ibge, ..."                       (footer note, blank)
#
# The fix removes those three rows outright (so every remaining row keeps
# its data), and renames the "unnamed: 1_level_1" column header to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank "situação do domicílio" row (original row 5).
$ws.Rows.Item(5).Delete()

# After the row-5 deletion everything shifted up by one, so the blank
# "grandes regiões e unidades da federação" row (originally row 8) is now
# row 7.
$ws.Rows.Item(7).Delete()

# After both deletions above, the blank footer "fonte: ibge, ..." row
# (originally row 41) is now row 39.
$ws.Rows.Item(39).Delete()

# Rename the stray pandas-style header label to "total".
$ws.Range("B2").Value() = "total"
